$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.982.87"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.819.61"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "310.49"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "0.4687"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").Value = "0.3669"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "0.07358"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "0.8740"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "20.30"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "1.826.62"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "5.419"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "0.07110"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "6.518"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "91.75"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "0.000008739"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "14.68"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "26.994.19"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "5.298"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").Value = "10.61"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "2.036.26"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").Value = "1.894"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "151.17"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "18.36"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "2.147"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").Value = "5.255"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "116.88"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "0.08890"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "0.7597"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").Value = "2.922"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D37").Value = "1.099"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "0.05310"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "0.01949"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "2.973"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "7.188"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "2.367"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("D43").Value = "0.5297"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").Value = "0.1654"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "8.473"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "0.4895"
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").Value = "1.665"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").Value = "103.32"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "0.06301"
$ws.Range("E51").Value = "  +0.10%  "
